# Justify all "regular text" paragraphs (title, abstract/body text, and the
# page-break-only paragraphs) while leaving the Heading 1 section titles
# untouched.
$wdAlignParagraphJustify = 3

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -ne "Heading 1") {
        $p.Format.Alignment = $wdAlignParagraphJustify
    }
}
